$wb = $excel.ActiveWorkbook

# --- About sheet ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B1").Value = "Arizona"
$wsAbout.Range("C1").NumberFormat = "mm-dd-yy"
$wsAbout.Range("C1").Value = 44319

# --- RTMF-passengers sheet ---
$wsPassengers = $wb.Worksheets.Item("RTMF-passengers")
$wsPassengers.Range("C2").Value = 0.15
$wsPassengers.Range("E2").Value = 0.05
$wsPassengers.Range("I2").Value = 0.8
$wsPassengers.Activate() | Out-Null
$wsPassengers.Range("I13").Select() | Out-Null

# --- RTMF-freight sheet becomes the active sheet/tab ---
$wsFreight = $wb.Worksheets.Item("RTMF-freight")
$wsFreight.Activate() | Out-Null
$wsFreight.Range("F4").Select() | Out-Null
